$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new week's task row
$ws.Range("B16").Value = "Francisco"

# Update existing task description to add a trailing period
$ws.Range("D9").Value = "Fixing issues with the pause menu, triggers and buttons."

$ws.Range("D16").Value = "Creating a new Player controller."

$ws.Range("D8").Select()
